$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column B ("Valor") needs to become column C, with a brand new
# "Variável" column inserted at B and a brand new "Colocação" column
# inserted after it (at D). Insert one column at a time so the existing
# data shifts by exactly one column each time.
$ws.Range("B1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("C1").Value = "Valor"
$ws.Range("D1").Value = "Colocação"

# "Variável" column: same label for every data row (2-10)
$ws.Range("B2:B10").Value = "Diferença 2021-2012"

# "Colocação" (ranking) column: only rows 2-8 get a value
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "19º"
